# virtualization-software-instance.xlsx edit
#
# Summary of the change being applied (per the commit "I am adding a bunch
# of files. I created a new simple database for virtual machine."):
#   - The VirtualBox "version" cell (B3) is updated to also include the
#     product name, matching the style already used for the VMware row.
#   - Column B is widened/auto-fit so the longer text is fully visible
#     (mirrors the user double-clicking the column border after editing).
#   - The active selection ends up on B3, the cell that was just edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the VirtualBox version string in B3 to be prefixed with the
#    product name, just like "VMware(r) Workstation 17 Pro ..." in B2.
$ws.Range("B3").Value = "VirtualBox Version 7.1.2 r164945 (Qt6.5.3)"

# 2. Widen column B to fit the new, longer text (best-fit-style column
#    resize). 44.66 is the input that this host's pixel-quantised
#    ColumnWidth setter rounds to the width closest to the real Excel
#    "best fit" result for this text.
$ws.Columns("B").ColumnWidth = 44.66

# 3. Leave the selection on the cell that was edited (B3), matching the
#    saved sheet view.
$ws.Range("B3").Select() | Out-Null
